$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 10-37 need the same bordered/centered/bold style that column A's
# existing cells (rows 2-9) use (style index 1 in the source file) -- copy
# that formatting down from A9 before writing any new values.
$ws.Range("A9").Copy()
$ws.Range("A10:A37").PasteSpecial(-4122)

# Column A holds short numeric-looking labels ("0", "1", "33", ...) that are
# stored as text in the workbook (t="inlineStr"), not as numbers. Force the
# whole label column to a text number format first so the values below are
# written as text instead of being auto-coerced to numbers.
$ws.Range("A2:A37").NumberFormat = "@"

# r, A (label, text), B (win rate), C (avg. game length)
$data = @(
  @(2, "0", 0, 33),
  @(3, "1", 100, 4),
  @(4, "2", 0, 24),
  @(5, "3", 100, 5),
  @(6, "4", 100, 2),
  @(7, "5", 100, 9),
  @(8, "7", 100, 11),
  @(9, "9", 0, 3),
  @(10, "11", 100, 22),
  @(11, "13", 100, 23),
  @(12, "15", 100, 11),
  @(13, "17", 100, 24),
  @(14, "19", 100, 21),
  @(15, "21", 0, 22),
  @(16, "23", 100, 5),
  @(17, "25", 0, 3),
  @(18, "27", 100, 5),
  @(19, "29", 0, 22),
  @(20, "31", 0, 23),
  @(21, "33", 100, 9),
  @(22, "35", 100, 19),
  @(23, "37", 100, 21),
  @(24, "39", 100, 14),
  @(25, "41", 0, 22),
  @(26, "43", 100, 15),
  @(27, "45", 100, 11),
  @(28, "46", 100, 11),
  @(29, "47", 100, 3),
  @(30, "48", 0, 21),
  @(31, "49", 100, 3),
  @(32, "50", 0, 2),
  @(33, "51", 100, 5),
  @(34, "52", 100, 5),
  @(35, "53", 100, 22),
  @(36, "54", 100, 4),
  @(37, "55", 0, 21)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
}
